# Scheduled market-price refresh for the Lich_Profits workbook.
# For each leve row below, currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the
# resulting LeveProfit(NQ/HQ) columns (H-N) are updated to the latest values.
# Some profit cells are removed/added entirely where the source no longer
# emits a value for that column (matches the authoritative OOXML diff).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart
$ws.Range("H19").Value = 1147.7858
$ws.Range("I19").Value = 448.0625
$ws.Range("J19").Value = 2080.75
$ws.Range("K19").Value = 448.0625
$ws.Range("L19").Value = 2080.75
$ws.Range("M19").Value = -273.0625
$ws.Range("N19").Value = -2430.75

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 8204.120000000001
$ws.Range("I137").Value = 10116.417
$ws.Range("K137").Value = 30349.251
$ws.Range("M137").Value = -27799.251

# Row 138: All-night Crafting
$ws.Range("H138").Value = 1922.55
$ws.Range("J138").Value = 2201.3384
$ws.Range("L138").Value = 6604.0152
$ws.Range("N138").Value = -16884.0152

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 3250.7837
$ws.Range("I2").Value = 4009.6553
$ws.Range("J2").Value = 499.875
$ws.Range("K2").Value = 4009.6553
$ws.Range("L2").Value = 499.875
$ws.Range("M2").Value = -3896.6553
$ws.Range("N2").Value = -725.875

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 8776.200000000001
$ws.Range("I61").Value = 7081.8335
$ws.Range("K61").Value = 7081.8335
$ws.Range("M61").Value = -6869.8335

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 2524.2273
$ws.Range("I63").Value = 2531.3333
$ws.Range("K63").Value = 2531.3333
$ws.Range("M63").Value = -1845.3333

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2524.2273
$ws.Range("I66").Value = 2531.3333
$ws.Range("K66").Value = 12656.6665
$ws.Range("M66").Value = -9224.666499999999

# Row 76: Sometimes the South Wins
$ws.Range("H76").Value = 5000
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5676

# Row 79: The Thriller of Autumn (L)
$ws.Range("H79").Value = 5000
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7340

# Row 116: No Scope
$ws.Range("H116").Value = 3250.7837
$ws.Range("I116").Value = 4009.6553
$ws.Range("J116").Value = 499.875
$ws.Range("K116").Value = 4009.6553
$ws.Range("L116").Value = 499.875
$ws.Range("M116").Value = -1715.6553
$ws.Range("N116").Value = -5087.875

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 8906.317999999999
$ws.Range("I122").Value = 5612.353
$ws.Range("K122").Value = 16837.059
$ws.Range("M122").Value = -14387.059

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 5162.081
$ws.Range("I132").Value = 3427.0605
$ws.Range("J132").Value = 19476
$ws.Range("K132").Value = 10281.1815
$ws.Range("L132").Value = 58428
$ws.Range("M132").Value = -7751.181500000001
$ws.Range("N132").Value = -63488

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 8776.200000000001
$ws.Range("I136").Value = 7081.8335
$ws.Range("K136").Value = 21245.5005
$ws.Range("M136").Value = -18695.5005

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 3250.7837
$ws.Range("I3").Value = 4009.6553
$ws.Range("J3").Value = 499.875
$ws.Range("K3").Value = 4009.6553
$ws.Range("L3").Value = 499.875
$ws.Range("M3").Value = -3895.6553
$ws.Range("N3").Value = -727.875

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 1905.6923
$ws.Range("I105").Value = 1814.5
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 1814.5
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -67.5
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 461.93103
$ws.Range("J22").Value = 846
$ws.Range("L22").Value = 846
$ws.Range("N22").Value = -1546

# Row 31: Wall Not Found
$ws.Range("H31").Value = 28541.064
$ws.Range("J31").Value = 38816.863
$ws.Range("L31").Value = 38816.863
$ws.Range("N31").Value = -39406.863

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 28541.064
$ws.Range("J34").Value = 38816.863
$ws.Range("L34").Value = 38816.863
$ws.Range("N34").Value = -39220.863

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3717.3572
$ws.Range("I58").Value = 4658
$ws.Range("J58").Value = 3011.875
$ws.Range("K58").Value = 4658
$ws.Range("L58").Value = 3011.875
$ws.Range("M58").Value = -4455
$ws.Range("N58").Value = -3417.875

# Row 99: O Pine
$ws.Range("H99").Value = 3046.1538
$ws.Range("J99").Value = 2979.1667
$ws.Range("L99").Value = 2979.1667
$ws.Range("N99").Value = -5975.1667

# Row 126: A Better Conductor
$ws.Range("H126").Value = 3046.1538
$ws.Range("J126").Value = 2979.1667
$ws.Range("L126").Value = 8937.500100000001
$ws.Range("N126").Value = -13877.5001

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 4031.348
$ws.Range("I132").Value = 3841.75
$ws.Range("J132").Value = 5295.3335
$ws.Range("K132").Value = 11525.25
$ws.Range("L132").Value = 15886.0005
$ws.Range("M132").Value = -8995.25
$ws.Range("N132").Value = -20946.0005

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 5785.48
$ws.Range("I134").Value = 6544
$ws.Range("K134").Value = 19632
$ws.Range("M134").Value = -17097

# Row 136: Turali Quality
$ws.Range("H136").Value = 3717.3572
$ws.Range("I136").Value = 4658
$ws.Range("J136").Value = 3011.875
$ws.Range("K136").Value = 13974
$ws.Range("L136").Value = 9035.625
$ws.Range("M136").Value = -11424
$ws.Range("N136").Value = -14135.625

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 105250.17
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 105250.17
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 105250.17
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -115610.17

$ws = $wb.Worksheets.Item("CUL")
# Row 59: Comfort Me with Mushrooms
$ws.Range("H59").Value = 1200
$ws.Range("I59").Value = 1200
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 3600
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -3060
$ws.Range("N59").ClearContents()

# Row 118: Teetotally
$ws.Range("H118").Value = 93.5
$ws.Range("I118").Value = 93.5
$ws.Range("K118").Value = 280.5
$ws.Range("M118").Value = 962.5

# Row 119: Super Dark Times
$ws.Range("H119").Value = 2197.5
$ws.Range("I119").Value = 2197.5
$ws.Range("K119").Value = 6592.5
$ws.Range("M119").Value = -1754.5

# Row 120: A Happy End
$ws.Range("H120").Value = 2950
$ws.Range("I120").Value = 2950
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 8850
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -4012
$ws.Range("N120").ClearContents()

# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 20837042
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 20837042
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 62511126
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -62513746

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 4121.722
$ws.Range("I80").Value = 3662.182
$ws.Range("K80").Value = 3662.182
$ws.Range("M80").Value = -2664.182

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4121.722
$ws.Range("I83").Value = 3662.182
$ws.Range("K83").Value = 18310.91
$ws.Range("M83").Value = -13318.91

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 1679.6129
$ws.Range("I102").Value = 1789.9286
$ws.Range("J102").Value = 650
$ws.Range("K102").Value = 1789.9286
$ws.Range("L102").Value = 650
$ws.Range("M102").Value = -167.9286
$ws.Range("N102").Value = -3894

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 9792.478999999999
$ws.Range("J126").Value = 8732.929
$ws.Range("L126").Value = 26198.787
$ws.Range("N126").Value = -31138.787

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 7640.8184
$ws.Range("J7").Value = 7066
$ws.Range("L7").Value = 7066
$ws.Range("N7").Value = -7290

# Row 40: Best Served Toad
$ws.Range("H40").Value = 4721.6523
$ws.Range("I40").Value = 4721.6523
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4721.6523
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4585.6523
$ws.Range("N40").ClearContents()

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 2642.4595
$ws.Range("J68").Value = 2691.6667
$ws.Range("L68").Value = 2691.6667
$ws.Range("N68").Value = -4189.6667

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 2642.4595
$ws.Range("J71").Value = 2691.6667
$ws.Range("L71").Value = 13458.3335
$ws.Range("N71").Value = -20946.3335

# Row 96: Off the Cuff
$ws.Range("H96").Value = 49997
$ws.Range("J96").Value = 49997
$ws.Range("L96").Value = 49997
$ws.Range("N96").Value = -55489

# Row 97: Looking for Glove
$ws.Range("H97").Value = 27500
$ws.Range("I97").Value = 25000
$ws.Range("J97").Value = 30000
$ws.Range("K97").Value = 25000
$ws.Range("L97").Value = 30000
$ws.Range("M97").Value = -24009
$ws.Range("N97").Value = -31982

# Row 103: Security Breeches
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 122: Hell on Leather
$ws.Range("H122").Value = 2797
$ws.Range("I122").Value = 2797
$ws.Range("K122").Value = 8391
$ws.Range("M122").Value = -5941

# Row 126: Battered Books
$ws.Range("H126").Value = 7640.8184
$ws.Range("J126").Value = 7066
$ws.Range("L126").Value = 21198
$ws.Range("N126").Value = -26138

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 3001.6843
$ws.Range("J136").Value = 4116.25
$ws.Range("L136").Value = 12348.75
$ws.Range("N136").Value = -17448.75

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax
$ws.Range("H107").Value = 409.3125
$ws.Range("I107").Value = 441.625
$ws.Range("J107").Value = 377
$ws.Range("K107").Value = 1324.875
$ws.Range("L107").Value = 1131
$ws.Range("M107").Value = 595.125
$ws.Range("N107").Value = -4971

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 3735.8333
$ws.Range("I126").Value = 2012.25
$ws.Range("K126").Value = 6036.75
$ws.Range("M126").Value = -3566.75
